# pontos notáveis - incremento na tabela de ranking
# Multiply the "percentage" columns E and F (rows 2-7) by 100,
# turning the stored fraction (0.xx) into a percentage number (xx.xx),
# while leaving the cell's number format (0.00%) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 7; $row++) {
    foreach ($col in @("E", "F")) {
        $cell = $ws.Range("$col$row")
        $cell.Value2 = $cell.Value2 * 100
    }
}
